$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source file names in column A were originally exported with a ".xmi"
# extension; correct them to use the ".tsv" extension instead (file endings
# fix), leaving every other cell/value untouched.
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Value2
    if ($oldVal -like "*.xmi") {
        $newVal = $oldVal -replace '\.xmi$', '.tsv'
        $cell.Value = $newVal
    }
}

# Match the author's final selection state: the whole of column A selected.
$ws.Range("A1:A1048576").Select()
